$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values: row 12 (Qps) and row 13 (Qs) column C
$ws.Range("C12").Value = 0.9
$ws.Range("C13").Value = 0.9

# Update the active cell selection shown in the sheetView
$ws.Range("G10").Select()
